$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110 (shifts existing rows 110-221 down to 111-222)
$ws.Rows("110:110").Insert()

# Populate the newly inserted row 110 with its data
$ws.Cells.Item(110, 1).Value  = 3
$ws.Cells.Item(110, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(110, 3).Value  = "Coquimbo"
$ws.Cells.Item(110, 4).Value  = 44874
$ws.Cells.Item(110, 5).Value  = 5
$ws.Cells.Item(110, 6).Value  = 100112026
$ws.Cells.Item(110, 7).Value  = "Haba"
$ws.Cells.Item(110, 8).Value  = "Sin especificar"
$ws.Cells.Item(110, 9).Value  = "Primera"
$ws.Cells.Item(110, 10).Value = 55
$ws.Cells.Item(110, 11).Value = 9000
$ws.Cells.Item(110, 12).Value = 9000
$ws.Cells.Item(110, 13).Value = 9000
$ws.Cells.Item(110, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(110, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(110, 16).Value = 360
$ws.Cells.Item(110, 17).Value = 25
$ws.Cells.Item(110, 18).Value = "Hortaliza"
